# Continue the data frame: append the newly-collected observation (row 38)
# to the Wetlands_Baro_complete sheet, matching the formatting already
# used for the surrounding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = -4122

# --- New data row ---
$newRow = 38

$ws.Range("A" + $newRow).Value = 9

$ws.Range("B" + $newRow).Value = 44750
$ws.Range("B3").Copy()
$ws.Range("B" + $newRow).PasteSpecial($xlPasteFormats)

$ws.Range("C" + $newRow).Value = 0.47847222222222219
$ws.Range("C3").Copy()
$ws.Range("C" + $newRow).PasteSpecial($xlPasteFormats)

$ws.Range("D" + $newRow).Value = 0.48541666666666666
$ws.Range("D3").Copy()
$ws.Range("D" + $newRow).PasteSpecial($xlPasteFormats)

$ws.Range("E" + $newRow).Value = 63.806420000000003
$ws.Range("F" + $newRow).Value = 8.6296610000000005

# The air-press/air-temp readings for this entry get a small distinguishing
# font treatment (9pt Lucida Sans, black) like other highlighted rows.
$ws.Range("E" + $newRow).Font.Color = 0
$ws.Range("E" + $newRow).Font.Size = 9
$ws.Range("E" + $newRow).Font.Name = "Lucida Sans"

$ws.Range("F" + $newRow).Font.Color = 0
$ws.Range("F" + $newRow).Font.Size = 9
$ws.Range("F" + $newRow).Font.Name = "Lucida Sans"

$ws.Application.CutCopyMode = $false

$ws.Range("C39").Select() | Out-Null
